$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'" + '63.569.18'
$ws.Cells.Item(2, 5).Value = "'" + '  -2.94%  '

$ws.Cells.Item(3, 4).Value = "'" + '3.274.37'
$ws.Cells.Item(3, 5).Value = "'" + '  -2.00%  '

$ws.Cells.Item(4, 4).Value = "'" + '0.997'
$ws.Cells.Item(4, 5).Value = "'" + '  -0.20%  '

$ws.Cells.Item(5, 4).Value = "'" + '522.48'
$ws.Cells.Item(5, 5).Value = "'" + '  -2.50%  '

$ws.Cells.Item(6, 4).Value = "'" + '170.03'
$ws.Cells.Item(6, 5).Value = "'" + '  -8.57%  '

$ws.Cells.Item(7, 4).Value = "'" + '0.581'
$ws.Cells.Item(7, 5).Value = "'" + '  -4.72%  '

$ws.Cells.Item(8, 2).Value = 'LidoStakedEther'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Cells.Item(8, 4).Value = "'" + '3.264.98'
$ws.Cells.Item(8, 5).Value = "'" + '  -2.17%  '

$ws.Cells.Item(9, 2).Value = 'USDC'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Cells.Item(9, 4).Value = "'" + '0.998'
$ws.Cells.Item(9, 5).Value = "'" + '  -0.23%  '

$ws.Cells.Item(10, 4).Value = "'" + '0.596'
$ws.Cells.Item(10, 5).Value = "'" + '  -5.03%  '

$ws.Cells.Item(11, 4).Value = "'" + '52.02'
$ws.Cells.Item(11, 5).Value = "'" + '  -13.26%  '

$ws.Cells.Item(12, 4).Value = "'" + '0.131'
$ws.Cells.Item(12, 5).Value = "'" + '  -3.16%  '

$ws.Cells.Item(13, 4).Value = "'" + '0.0000255'
$ws.Cells.Item(13, 5).Value = "'" + '  -4.69%  '

$ws.Cells.Item(14, 4).Value = "'" + '8.82'
$ws.Cells.Item(14, 5).Value = "'" + '  -4.41%  '

$ws.Cells.Item(15, 4).Value = "'" + '3.756.75'
$ws.Cells.Item(15, 5).Value = "'" + '  -2.78%  '

$ws.Cells.Item(16, 2).Value = 'TRON'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(16, 4).Value = "'" + '0.116'
$ws.Cells.Item(16, 5).Value = "'" + '  -2.11%  '

$ws.Cells.Item(17, 2).Value = 'WrappedEther'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(17, 4).Value = "'" + '3.243.73'
$ws.Cells.Item(17, 5).Value = "'" + '  -2.82%  '

$ws.Cells.Item(18, 4).Value = "'" + '63.158.30'
$ws.Cells.Item(18, 5).Value = "'" + '  -3.04%  '

$ws.Cells.Item(19, 4).Value = "'" + '17.21'
$ws.Cells.Item(19, 5).Value = "'" + '  -3.85%  '

$ws.Cells.Item(20, 4).Value = "'" + '11.07'
$ws.Cells.Item(20, 5).Value = "'" + '  -1.99%  '

$ws.Cells.Item(21, 4).Value = "'" + '0.947'
$ws.Cells.Item(21, 5).Value = "'" + '  -2.46%  '

$ws.Cells.Item(22, 4).Value = "'" + '374.63'
$ws.Cells.Item(22, 5).Value = "'" + '  -1.20%  '

$ws.Cells.Item(23, 4).Value = "'" + '4.14'
$ws.Cells.Item(23, 5).Value = "'" + '  +6.17%  '

$ws.Cells.Item(24, 4).Value = "'" + '80.66'
$ws.Cells.Item(24, 5).Value = "'" + '  -1.40%  '

$ws.Cells.Item(25, 4).Value = "'" + '11.05'
$ws.Cells.Item(25, 5).Value = "'" + '  -3.36%  '

$ws.Cells.Item(26, 4).Value = "'" + '3.64'
$ws.Cells.Item(26, 5).Value = "'" + '  -5.97%  '

$ws.Cells.Item(27, 5).Value = "'" + '  +2.10%  '

$ws.Cells.Item(28, 4).Value = "'" + '2.67'
$ws.Cells.Item(28, 5).Value = "'" + '  -1.93%  '

$ws.Cells.Item(29, 4).Value = "'" + '11.08'
$ws.Cells.Item(29, 5).Value = "'" + '  -5.29%  '

$ws.Cells.Item(30, 4).Value = "'" + '8.00'
$ws.Cells.Item(30, 5).Value = "'" + '  -6.55%  '

$ws.Cells.Item(31, 4).Value = "'" + '28.38'
$ws.Cells.Item(31, 5).Value = "'" + '  -3.24%  '

$ws.Cells.Item(32, 4).Value = "'" + '619.93'
$ws.Cells.Item(32, 5).Value = "'" + '  -5.03%  '

$ws.Cells.Item(33, 4).Value = "'" + '6.50'
$ws.Cells.Item(33, 5).Value = "'" + '  -5.75%  '

$ws.Cells.Item(34, 4).Value = "'" + '11.07'
$ws.Cells.Item(34, 5).Value = "'" + '  -3.17%  '

$ws.Cells.Item(35, 4).Value = "'" + '0.104'
$ws.Cells.Item(35, 5).Value = "'" + '  -3.11%  '

$ws.Cells.Item(36, 4).Value = "'" + '56.11'
$ws.Cells.Item(36, 5).Value = "'" + '  -6.47%  '

$ws.Cells.Item(37, 4).Value = "'" + '1.00'
$ws.Cells.Item(37, 5).Value = "'" + '  +0.13%  '

$ws.Cells.Item(38, 4).Value = "'" + '35.43'
$ws.Cells.Item(38, 5).Value = "'" + '  -4.57%  '

$ws.Cells.Item(39, 4).Value = "'" + '0.371'
$ws.Cells.Item(39, 5).Value = "'" + '  -7.13%  '

$ws.Cells.Item(40, 2).Value = 'PEPE'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(40, 4).Value = "'" + '0.0₃0730'
$ws.Cells.Item(40, 5).Value = "'" + '  -0.87%  '

$ws.Cells.Item(41, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(41, 4).Value = "'" + '0.996'
$ws.Cells.Item(41, 5).Value = "'" + '  -0.07%  '

$ws.Cells.Item(42, 2).Value = 'Stacks'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(42, 4).Value = "'" + '3.11'
$ws.Cells.Item(42, 5).Value = "'" + '  +6.53%  '

$ws.Cells.Item(43, 2).Value = 'Fetch.AI'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(43, 4).Value = "'" + '2.58'
$ws.Cells.Item(43, 5).Value = "'" + '  +0.97%  '

$ws.Cells.Item(44, 2).Value = 'Kaspa'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(44, 4).Value = "'" + '0.123'
$ws.Cells.Item(44, 5).Value = "'" + '  -4.88%  '

$ws.Cells.Item(45, 4).Value = "'" + '2.823.49'
$ws.Cells.Item(45, 5).Value = "'" + '  -3.53%  '

$ws.Cells.Item(46, 4).Value = "'" + '2.66'
$ws.Cells.Item(46, 5).Value = "'" + '  -1.14%  '

$ws.Cells.Item(47, 4).Value = "'" + '0.0391'
$ws.Cells.Item(47, 5).Value = "'" + '  -3.78%  '

$ws.Cells.Item(48, 2).Value = 'ThetaToken'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Cells.Item(48, 4).Value = "'" + '2.57'
$ws.Cells.Item(48, 5).Value = "'" + '  -6.08%  '

$ws.Cells.Item(49, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(49, 4).Value = "'" + '2.96'
$ws.Cells.Item(49, 5).Value = "'" + '  -1.29%  '

$ws.Cells.Item(50, 4).Value = "'" + '136.48'
$ws.Cells.Item(50, 5).Value = "'" + '  +0.79%  '

$ws.Cells.Item(51, 4).Value = "'" + '0.123'
